$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "0.013±0.002"
$ws.Range("C2").Value = "0.205±0.009"

$ws.Range("B3").Value = "0.075±0.004"
$ws.Range("C3").Value = "0.157±0.026"

$ws.Range("B4").Value = "0.669±0.058"
$ws.Range("C4").Value = "0.175±0.016"

$ws.Range("B5").Value = "0.968±0.007"
$ws.Range("C5").Value = "0.401±0.016"

$ws.Range("B6").Value = "0.892±0.015"
$ws.Range("C6").Value = "0.604±0.050"

$ws.Range("B7").Value = "0.729±0.060"
$ws.Range("C7").Value = "0.104±0.012"

$ws.Range("B8").Value = "0.005±0.001"
$ws.Range("C8").Value = "0.207±0.007"
